# Error Calculations and Plots
# Two data rows ("RM 232" and "SC 92") are removed from the missing-data
# table, shifting the subsequent rows up. In addition, the set of
# "missing" values in column F is changed: the row that is now "SC 5"
# (row 26) gets its previously-missing F value filled in with 17.38,
# while the row that is now "SC 101" (row 27) has its F value cleared
# out (becoming the new missing value).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete "SC 92" (originally row 28) first so row numbers for the
# still-to-be-deleted "RM 232" row (26) remain valid.
$ws.Rows.Item(28).Delete()
$ws.Rows.Item(26).Delete()

# After the deletions, "SC 5" is now row 26 and "SC 101" is row 27.
$ws.Range("F26").Value = 17.38
$ws.Range("F27").ClearContents()
